# live_trading_results.xlsx update
# Trade #49 (global trade #107 / momentum trade #78) closed at 2026-02-18 00:20:24
# and a brand-new momentum trade #107 (global) opened at 2026-02-18 00:20:18.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.45   # Current Capital
$summary.Range("B4").Value = 0.55      # Total P&L $
$summary.Range("B5").Value = 0.14      # Total P&L %
$summary.Range("B6").Value = 77        # Total Trades
$summary.Range("B8").Value = 31        # Losing Trades
$summary.Range("B9").Value = 49.35     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - "momentum" strategy row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.65000000000001   # Capital
$status.Range("D11").Value = 10                  # Trades
$status.Range("E11").Value = -0.35               # P&L $
$status.Range("F11").Value = -0.35               # P&L %
$status.Range("G11").Value = 10                  # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Row 79 (Trade # 78) transitions from OPEN to CLOSED
$allTrades.Range("G79").Value = 0.64
$allTrades.Range("H79").Value = "CLOSED"
$allTrades.Range("I79").Value = -4.4776
$allTrades.Range("J79").Value = -0.03
$allTrades.Range("K79").Value = 99.65000000000001
$allTrades.Range("L79").Value = "early_exit"
$allTrades.Range("M79").Value = 0.13

# New row 108 (Trade # 107) - newly opened momentum trade
$allTrades.Range("A108").Value = 107
$allTrades.Range("B108").NumberFormat = "@"
$allTrades.Range("B108").Value = "2026-02-18"
$allTrades.Range("B108").Style = "Normal"
$allTrades.Range("C108").Value = "00:20:18"
$allTrades.Range("D108").Value = "momentum"
$allTrades.Range("E108").Value = "DOWN"
$allTrades.Range("F108").Value = 0.67
$allTrades.Range("G108").Value = ""
$allTrades.Range("H108").Value = "OPEN"
$allTrades.Range("I108").Value = 0
$allTrades.Range("J108").Value = 0
$allTrades.Range("K108").Value = 99.6787371310913
$allTrades.Range("L108").Value = ""
$allTrades.Range("M108").Value = 0
$allTrades.Range("N108").Value = 0
$allTrades.Range("O108").Value = 0
$allTrades.Range("P108").Value = 0.9
$allTrades.Range("Q108").Value = "Downward momentum: -3.810% over 10 samples"

# ---------------------------------------------------------------------------
# momentum sheet
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

# Row 12 (Trade # 78) transitions from OPEN to CLOSED
$momentum.Range("G12").Value = 0.64
$momentum.Range("H12").Value = "CLOSED"
$momentum.Range("I12").Value = -4.4776
$momentum.Range("J12").Value = -0.03
$momentum.Range("K12").Value = 99.65000000000001
$momentum.Range("P12").Value = "early_exit"
$momentum.Range("Q12").Value = 0.13

# New row 26 (Trade # 107) - newly opened momentum trade
$momentum.Range("A26").Value = 107
$momentum.Range("B26").NumberFormat = "@"
$momentum.Range("B26").Value = "2026-02-18"
$momentum.Range("B26").Style = "Normal"
$momentum.Range("C26").Value = "00:20:18"
$momentum.Range("D26").Value = "momentum"
$momentum.Range("E26").Value = "DOWN"
$momentum.Range("F26").Value = 0.67
$momentum.Range("G26").Value = ""
$momentum.Range("H26").Value = "OPEN"
$momentum.Range("I26").Value = 0
$momentum.Range("J26").Value = 0
$momentum.Range("K26").Value = 99.6787371310913
$momentum.Range("L26").Value = 0
$momentum.Range("M26").Value = 0
$momentum.Range("N26").Value = 0.9
$momentum.Range("O26").Value = "Downward momentum: -3.810% over 10 samples"
$momentum.Range("P26").Value = ""
$momentum.Range("Q26").Value = 0
